$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update status + metrics
$ws.Range("B2").Value = "complete"
$ws.Range("J2").Value = "[81, 91]"
$ws.Range("K2").Value = "[0.086, 0.086]"
$ws.Range("L2").Value = "[0, 5]"
$ws.Range("M2").Value = "[0, 2.3000000000000007]"
$ws.Range("N2").Value = "[0, 0.1]"

# Rows 3-10: clear the result columns (B through O), keep column A (index) intact
$ws.Range("B3:O10").ClearContents()
